# Updates the cryptos price/volume table to the latest scrape.
# All price (D) and volume-change (E) cells are stored as plain text in the
# source sheet (e.g. "1.006", "27.118.06", "  -1.48%  "), so each write
# forces text formatting before assigning the value (otherwise Excel's
# automatic type inference would turn plain-decimal-looking strings such as
# "1.006" into the number 1.006) and restores the cell to the workbook's
# default "Normal" style afterwards so no stray formatting is introduced.
# Rows 44/45 additionally swap the Frax / EnergySwap coin rows (name, link,
# price and volume all change together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.118.06' },
    @{ Cell = 'E2'; Value = '  -1.48%  ' },
    @{ Cell = 'D3'; Value = '1.780.37' },
    @{ Cell = 'E3'; Value = '  -2.10%  ' },
    @{ Cell = 'D4'; Value = '1.006' },
    @{ Cell = 'E4'; Value = '  +0.24%  ' },
    @{ Cell = 'D5'; Value = '336.45' },
    @{ Cell = 'E5'; Value = '  -2.57%  ' },
    @{ Cell = 'D6'; Value = '1.003' },
    @{ Cell = 'E6'; Value = '  +0.27%  ' },
    @{ Cell = 'D7'; Value = '0.3812' },
    @{ Cell = 'E7'; Value = '  -0.33%  ' },
    @{ Cell = 'D8'; Value = '0.3413' },
    @{ Cell = 'E8'; Value = '  -3.12%  ' },
    @{ Cell = 'D9'; Value = '48.07' },
    @{ Cell = 'E9'; Value = '  -3.52%  ' },
    @{ Cell = 'D10'; Value = '1.186' },
    @{ Cell = 'E10'; Value = '  -4.17%  ' },
    @{ Cell = 'D11'; Value = '0.07432' },
    @{ Cell = 'E11'; Value = '  -4.67%  ' },
    @{ Cell = 'D12'; Value = '1.003' },
    @{ Cell = 'E12'; Value = '  +0.08%  ' },
    @{ Cell = 'D13'; Value = '21.61' },
    @{ Cell = 'E13'; Value = '  -3.21%  ' },
    @{ Cell = 'D14'; Value = '6.412' },
    @{ Cell = 'E14'; Value = '  -3.24%  ' },
    @{ Cell = 'D15'; Value = '1.779.12' },
    @{ Cell = 'E15'; Value = '  -1.92%  ' },
    @{ Cell = 'D16'; Value = '7.049' },
    @{ Cell = 'E16'; Value = '  -2.50%  ' },
    @{ Cell = 'D17'; Value = '0.00001083' },
    @{ Cell = 'E17'; Value = '  -3.97%  ' },
    @{ Cell = 'D18'; Value = '0.06663' },
    @{ Cell = 'E18'; Value = '  -0.94%  ' },
    @{ Cell = 'D19'; Value = '83.23' },
    @{ Cell = 'E19'; Value = '  -3.63%  ' },
    @{ Cell = 'E20'; Value = '  +0.31%  ' },
    @{ Cell = 'D21'; Value = '6.546' },
    @{ Cell = 'E21'; Value = '  +0.10%  ' },
    @{ Cell = 'D22'; Value = '17.26' },
    @{ Cell = 'E22'; Value = '  -2.77%  ' },
    @{ Cell = 'D23'; Value = '27.116.80' },
    @{ Cell = 'E23'; Value = '  -1.46%  ' },
    @{ Cell = 'D24'; Value = '12.21' },
    @{ Cell = 'E24'; Value = '  -7.45%  ' },
    @{ Cell = 'D25'; Value = '2.387' },
    @{ Cell = 'E25'; Value = '  -3.07%  ' },
    @{ Cell = 'D26'; Value = '2.504' },
    @{ Cell = 'E26'; Value = '  -6.65%  ' },
    @{ Cell = 'D27'; Value = '1.464' },
    @{ Cell = 'E27'; Value = '  -2.15%  ' },
    @{ Cell = 'D28'; Value = '21.06' },
    @{ Cell = 'E28'; Value = '  -5.00%  ' },
    @{ Cell = 'D29'; Value = '154.64' },
    @{ Cell = 'E29'; Value = '  +0.77%  ' },
    @{ Cell = 'D30'; Value = '1.978.85' },
    @{ Cell = 'E30'; Value = '  -1.92%  ' },
    @{ Cell = 'D31'; Value = '133.99' },
    @{ Cell = 'E31'; Value = '  -1.87%  ' },
    @{ Cell = 'D32'; Value = '3.993' },
    @{ Cell = 'E32'; Value = '  -2.21%  ' },
    @{ Cell = 'D33'; Value = '6.006' },
    @{ Cell = 'E33'; Value = '  -5.94%  ' },
    @{ Cell = 'D34'; Value = '0.08652' },
    @{ Cell = 'E34'; Value = '  -1.58%  ' },
    @{ Cell = 'D35'; Value = '13.04' },
    @{ Cell = 'E35'; Value = '  -7.39%  ' },
    @{ Cell = 'D36'; Value = '1.618' },
    @{ Cell = 'E36'; Value = '  -5.18%  ' },
    @{ Cell = 'D37'; Value = '5.379' },
    @{ Cell = 'E37'; Value = '  -4.75%  ' },
    @{ Cell = 'D38'; Value = '0.6809' },
    @{ Cell = 'E38'; Value = '  -4.04%  ' },
    @{ Cell = 'D39'; Value = '0.06277' },
    @{ Cell = 'E39'; Value = '  -4.23%  ' },
    @{ Cell = 'D40'; Value = '0.2171' },
    @{ Cell = 'E40'; Value = '  -4.85%  ' },
    @{ Cell = 'D41'; Value = '0.02317' },
    @{ Cell = 'E41'; Value = '  -4.49%  ' },
    @{ Cell = 'D42'; Value = '8.581' },
    @{ Cell = 'E42'; Value = '  -4.96%  ' },
    @{ Cell = 'D43'; Value = '1.225' },
    @{ Cell = 'E43'; Value = '  -4.39%  ' },
    @{ Cell = 'B44'; Value = 'Frax' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax' },
    @{ Cell = 'D44'; Value = '1.003' },
    @{ Cell = 'E44'; Value = '  +0.34%  ' },
    @{ Cell = 'B45'; Value = 'EnergySwap' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D45'; Value = '14.15' },
    @{ Cell = 'E45'; Value = '  -4.16%  ' },
    @{ Cell = 'D46'; Value = '0.6397' },
    @{ Cell = 'E46'; Value = '  -3.64%  ' },
    @{ Cell = 'D47'; Value = '3.855' },
    @{ Cell = 'E47'; Value = '  -2.98%  ' },
    @{ Cell = 'D48'; Value = '2.125' },
    @{ Cell = 'E48'; Value = '  -2.73%  ' },
    @{ Cell = 'D49'; Value = '130.91' },
    @{ Cell = 'E49'; Value = '  -1.49%  ' },
    @{ Cell = 'D50'; Value = '0.07092' },
    @{ Cell = 'E50'; Value = '  -3.73%  ' },
    @{ Cell = 'D51'; Value = '78.54' },
    @{ Cell = 'E51'; Value = '  -2.77%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
